$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill rows 369-401 with new certificate records. Formatting is copied
#     from the existing alternating 9-row block (rows 360-368: style-12 rows use
#     row 360 as the template, style-11 rows use row 361), then the five
#     cell values (Certificate No, Name, Course, Date, Result) are written. ---
$ws.Range("A360:E360").Copy()
$ws.Range("A369:E369").PasteSpecial(-4122)
$ws.Cells.Item(369,1).Value = 'DSS1368'
$ws.Cells.Item(369,2).Value = 'RADWAN MOHAMED ABDELHAFEZ HASSAN'
$ws.Cells.Item(369,3).Value = 'Scaffold Competent Person'
$ws.Cells.Item(369,4).Value = 45517
$ws.Cells.Item(369,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A370:E370").PasteSpecial(-4122)
$ws.Cells.Item(370,1).Value = 'DSS1369'
$ws.Cells.Item(370,2).Value = 'MOHAMMED ALMOSTAFA MOHAMMED ELDAW'
$ws.Cells.Item(370,3).Value = '30 Hours Construction Safety & Health'
$ws.Cells.Item(370,4).Value = 45566
$ws.Cells.Item(370,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A371:E371").PasteSpecial(-4122)
$ws.Cells.Item(371,1).Value = 'DSS1370'
$ws.Cells.Item(371,2).Value = 'MOHAMMED ALMOSTAFA MOHAMMED ELDAW'
$ws.Cells.Item(371,3).Value = '30 Hours G. Industry Safety & Health'
$ws.Cells.Item(371,4).Value = 45566
$ws.Cells.Item(371,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A372:E372").PasteSpecial(-4122)
$ws.Cells.Item(372,1).Value = 'DSS1371'
$ws.Cells.Item(372,2).Value = 'MOHAMMED ALMOSTAFA MOHAMMED ELDAW'
$ws.Cells.Item(372,3).Value = 'Electrical Safety '
$ws.Cells.Item(372,4).Value = 45566
$ws.Cells.Item(372,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A373:E373").PasteSpecial(-4122)
$ws.Cells.Item(373,1).Value = 'DSS1372'
$ws.Cells.Item(373,2).Value = 'MOHAMMED ALMOSTAFA MOHAMMED ELDAW'
$ws.Cells.Item(373,3).Value = 'Fire Marshal'
$ws.Cells.Item(373,4).Value = 45566
$ws.Cells.Item(373,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A374:E374").PasteSpecial(-4122)
$ws.Cells.Item(374,1).Value = 'DSS1373'
$ws.Cells.Item(374,2).Value = 'MOHAMMED ALMOSTAFA MOHAMMED ELDAW'
$ws.Cells.Item(374,3).Value = 'Scaffold Competent Person'
$ws.Cells.Item(374,4).Value = 45566
$ws.Cells.Item(374,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A375:E375").PasteSpecial(-4122)
$ws.Cells.Item(375,1).Value = 'DSS1374'
$ws.Cells.Item(375,2).Value = 'MOHAMMED ALMOSTAFA MOHAMMED ELDAW'
$ws.Cells.Item(375,3).Value = 'Lifting & Rigging Competent Person'
$ws.Cells.Item(375,4).Value = 45566
$ws.Cells.Item(375,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A376:E376").PasteSpecial(-4122)
$ws.Cells.Item(376,1).Value = 'DSS1375'
$ws.Cells.Item(376,2).Value = 'MOHAMMED ALMOSTAFA MOHAMMED ELDAW'
$ws.Cells.Item(376,3).Value = 'Health & Safety Risk Assessment'
$ws.Cells.Item(376,4).Value = 45566
$ws.Cells.Item(376,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A377:E377").PasteSpecial(-4122)
$ws.Cells.Item(377,1).Value = 'DSS1376'
$ws.Cells.Item(377,2).Value = 'MOHAMMED ALMOSTAFA MOHAMMED ELDAW'
$ws.Cells.Item(377,3).Value = 'Safety Management System & PTW'
$ws.Cells.Item(377,4).Value = 45566
$ws.Cells.Item(377,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A378:E378").PasteSpecial(-4122)
$ws.Cells.Item(378,1).Value = 'DSS1377'
$ws.Cells.Item(378,2).Value = 'Mohamed Mahmoud Farag Mohamed'
$ws.Cells.Item(378,3).Value = '30 Hours Construction Safety & Health'
$ws.Cells.Item(378,4).Value = 45566
$ws.Cells.Item(378,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A379:E379").PasteSpecial(-4122)
$ws.Cells.Item(379,1).Value = 'DSS1378'
$ws.Cells.Item(379,2).Value = 'Mohamed Mahmoud Farag Mohamed'
$ws.Cells.Item(379,3).Value = '30 Hours G. Industry Safety & Health'
$ws.Cells.Item(379,4).Value = 45566
$ws.Cells.Item(379,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A380:E380").PasteSpecial(-4122)
$ws.Cells.Item(380,1).Value = 'DSS1379'
$ws.Cells.Item(380,2).Value = 'Mohamed Mahmoud Farag Mohamed'
$ws.Cells.Item(380,3).Value = 'Electrical Safety '
$ws.Cells.Item(380,4).Value = 45566
$ws.Cells.Item(380,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A381:E381").PasteSpecial(-4122)
$ws.Cells.Item(381,1).Value = 'DSS1380'
$ws.Cells.Item(381,2).Value = 'Mohamed Mahmoud Farag Mohamed'
$ws.Cells.Item(381,3).Value = 'Fire Marshal'
$ws.Cells.Item(381,4).Value = 45566
$ws.Cells.Item(381,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A382:E382").PasteSpecial(-4122)
$ws.Cells.Item(382,1).Value = 'DSS1381'
$ws.Cells.Item(382,2).Value = 'Mohamed Mahmoud Farag Mohamed'
$ws.Cells.Item(382,3).Value = 'Scaffold Competent Person'
$ws.Cells.Item(382,4).Value = 45566
$ws.Cells.Item(382,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A383:E383").PasteSpecial(-4122)
$ws.Cells.Item(383,1).Value = 'DSS1382'
$ws.Cells.Item(383,2).Value = 'Mohamed Mahmoud Farag Mohamed'
$ws.Cells.Item(383,3).Value = 'Lifting & Rigging Competent Person'
$ws.Cells.Item(383,4).Value = 45566
$ws.Cells.Item(383,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A384:E384").PasteSpecial(-4122)
$ws.Cells.Item(384,1).Value = 'DSS1383'
$ws.Cells.Item(384,2).Value = 'Mohamed Mahmoud Farag Mohamed'
$ws.Cells.Item(384,3).Value = 'Health & Safety Risk Assessment'
$ws.Cells.Item(384,4).Value = 45566
$ws.Cells.Item(384,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A385:E385").PasteSpecial(-4122)
$ws.Cells.Item(385,1).Value = 'DSS1384'
$ws.Cells.Item(385,2).Value = 'Mohamed Mahmoud Farag Mohamed'
$ws.Cells.Item(385,3).Value = 'Safety Management System & PTW'
$ws.Cells.Item(385,4).Value = 45566
$ws.Cells.Item(385,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A386:E386").PasteSpecial(-4122)
$ws.Cells.Item(386,1).Value = 'DSS1385'
$ws.Cells.Item(386,2).Value = 'Mohamed Abo Alhassan Mohamed Hassanen'
$ws.Cells.Item(386,3).Value = '30 Hours Construction Safety & Health'
$ws.Cells.Item(386,4).Value = 45566
$ws.Cells.Item(386,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A387:E387").PasteSpecial(-4122)
$ws.Cells.Item(387,1).Value = 'DSS1386'
$ws.Cells.Item(387,2).Value = 'Mohamed Abo Alhassan Mohamed Hassanen'
$ws.Cells.Item(387,3).Value = '30 Hours G. Industry Safety & Health'
$ws.Cells.Item(387,4).Value = 45566
$ws.Cells.Item(387,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A388:E388").PasteSpecial(-4122)
$ws.Cells.Item(388,1).Value = 'DSS1387'
$ws.Cells.Item(388,2).Value = 'Mohamed Abo Alhassan Mohamed Hassanen'
$ws.Cells.Item(388,3).Value = 'Electrical Safety '
$ws.Cells.Item(388,4).Value = 45566
$ws.Cells.Item(388,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A389:E389").PasteSpecial(-4122)
$ws.Cells.Item(389,1).Value = 'DSS1388'
$ws.Cells.Item(389,2).Value = 'Mohamed Abo Alhassan Mohamed Hassanen'
$ws.Cells.Item(389,3).Value = 'Fire Marshal'
$ws.Cells.Item(389,4).Value = 45566
$ws.Cells.Item(389,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A390:E390").PasteSpecial(-4122)
$ws.Cells.Item(390,1).Value = 'DSS1389'
$ws.Cells.Item(390,2).Value = 'Mohamed Abo Alhassan Mohamed Hassanen'
$ws.Cells.Item(390,3).Value = 'Scaffold Competent Person'
$ws.Cells.Item(390,4).Value = 45566
$ws.Cells.Item(390,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A391:E391").PasteSpecial(-4122)
$ws.Cells.Item(391,1).Value = 'DSS1390'
$ws.Cells.Item(391,2).Value = 'Mohamed Abo Alhassan Mohamed Hassanen'
$ws.Cells.Item(391,3).Value = 'Lifting & Rigging Competent Person'
$ws.Cells.Item(391,4).Value = 45566
$ws.Cells.Item(391,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A392:E392").PasteSpecial(-4122)
$ws.Cells.Item(392,1).Value = 'DSS1391'
$ws.Cells.Item(392,2).Value = 'Mohamed Abo Alhassan Mohamed Hassanen'
$ws.Cells.Item(392,3).Value = 'Health & Safety Risk Assessment'
$ws.Cells.Item(392,4).Value = 45566
$ws.Cells.Item(392,5).Value = 1

$ws.Range("A361:E361").Copy()
$ws.Range("A393:E393").PasteSpecial(-4122)
$ws.Cells.Item(393,1).Value = 'DSS1392'
$ws.Cells.Item(393,2).Value = 'Mohamed Abo Alhassan Mohamed Hassanen'
$ws.Cells.Item(393,3).Value = 'Safety Management System & PTW'
$ws.Cells.Item(393,4).Value = 45566
$ws.Cells.Item(393,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A394:E394").PasteSpecial(-4122)
$ws.Cells.Item(394,1).Value = 'DSS1393'
$ws.Cells.Item(394,2).Value = 'RAMADAN HUSSEIN MOHAMED AWAD'
$ws.Cells.Item(394,3).Value = '30 Hours Construction Safety & Health'
$ws.Cells.Item(394,4).Value = 45566
$ws.Cells.Item(394,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A395:E395").PasteSpecial(-4122)
$ws.Cells.Item(395,1).Value = 'DSS1394'
$ws.Cells.Item(395,2).Value = 'RAMADAN HUSSEIN MOHAMED AWAD'
$ws.Cells.Item(395,3).Value = '30 Hours G. Industry Safety & Health'
$ws.Cells.Item(395,4).Value = 45566
$ws.Cells.Item(395,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A396:E396").PasteSpecial(-4122)
$ws.Cells.Item(396,1).Value = 'DSS1395'
$ws.Cells.Item(396,2).Value = 'RAMADAN HUSSEIN MOHAMED AWAD'
$ws.Cells.Item(396,3).Value = 'Electrical Safety '
$ws.Cells.Item(396,4).Value = 45566
$ws.Cells.Item(396,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A397:E397").PasteSpecial(-4122)
$ws.Cells.Item(397,1).Value = 'DSS1396'
$ws.Cells.Item(397,2).Value = 'RAMADAN HUSSEIN MOHAMED AWAD'
$ws.Cells.Item(397,3).Value = 'Fire Marshal'
$ws.Cells.Item(397,4).Value = 45566
$ws.Cells.Item(397,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A398:E398").PasteSpecial(-4122)
$ws.Cells.Item(398,1).Value = 'DSS1397'
$ws.Cells.Item(398,2).Value = 'RAMADAN HUSSEIN MOHAMED AWAD'
$ws.Cells.Item(398,3).Value = 'Scaffold Competent Person'
$ws.Cells.Item(398,4).Value = 45566
$ws.Cells.Item(398,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A399:E399").PasteSpecial(-4122)
$ws.Cells.Item(399,1).Value = 'DSS1398'
$ws.Cells.Item(399,2).Value = 'RAMADAN HUSSEIN MOHAMED AWAD'
$ws.Cells.Item(399,3).Value = 'Lifting & Rigging Competent Person'
$ws.Cells.Item(399,4).Value = 45566
$ws.Cells.Item(399,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A400:E400").PasteSpecial(-4122)
$ws.Cells.Item(400,1).Value = 'DSS1399'
$ws.Cells.Item(400,2).Value = 'RAMADAN HUSSEIN MOHAMED AWAD'
$ws.Cells.Item(400,3).Value = 'Health & Safety Risk Assessment'
$ws.Cells.Item(400,4).Value = 45566
$ws.Cells.Item(400,5).Value = 1

$ws.Range("A360:E360").Copy()
$ws.Range("A401:E401").PasteSpecial(-4122)
$ws.Cells.Item(401,1).Value = 'DSS1400'
$ws.Cells.Item(401,2).Value = 'RAMADAN HUSSEIN MOHAMED AWAD'
$ws.Cells.Item(401,3).Value = 'Safety Management System & PTW'
$ws.Cells.Item(401,4).Value = 45566
$ws.Cells.Item(401,5).Value = 1

$excel.CutCopyMode = $false

# --- Remove one trailing blank row (row 695); the old row 696 (thick bottom
#     border row) shifts up to become the new last row, 695. ---
$ws.Rows.Item(695).Delete()

# --- Refresh AutoFilter so its range covers the new extent A1:E695 ---
if ($ws.AutoFilterMode) {
  $ws.AutoFilterMode = $false
}
$ws.Range("A1:E695").AutoFilter()

# --- Shrink the _FilterDatabase / Print_Area named ranges to match ---
$wb.Names.Item(1).RefersTo = '=Sheet1!$A$1:$E$695'
$wb.Names.Item(2).RefersTo = '=Sheet1!$A$1:$E$695'

# --- Update the window/view state: scroll so row 384 is at the top and select F391 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 384
$win.ScrollColumn = 1
$ws.Range("F391").Select()

Write-Host "Edit complete"